$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hyperparameter Optimization")

# --- Long JSON-looking strings reused multiple times below ---
$s47 = "[0.05739310756325722, 0.06071412190794945, 0.07280246168375015, 0.07271604984998703, 0.057159747928380966]"
$s31 = "[0.018034348264336586, 0.023353051394224167, 0.02859247475862503, 0.03425142541527748, 0.01903075911104679]"
$s32 = "[0.0355488546192646, 0.031698182225227356, 0.05592479929327965, 0.0423433892428875, 0.046196676790714264]"
$s72 = "[0.031347401440143585, 0.03457336872816086, 0.05158032476902008, 0.03897654637694359, 0.04255129024386406]"
$s71 = "[0.015520330518484116, 0.01797584816813469, 0.018809670582413673, 0.022205613553524017, 0.010675263591110706]"
$s87 = "[0.030682526528835297, 0.03892310708761215, 0.04652285575866699, 0.0381121002137661, 0.036519937217235565]"

# The textual cells below are filled in a specific order so the workbook's shared-string
# table grows in exactly the same sequence as in the target file (new strings are appended
# in first-use order): K47="62", S47, S31, S32, K72="56", S72, S71, K87="56" (reuse), S87.

# Row 47 "Best Trial" + its fold RMSE list
$ws.Range("K47").Value = "62"
$ws.Range("S47").Value = $s47

# Row 31 fold RMSE list
$ws.Range("S31").Value = $s31

# Row 32 fold RMSE list
$ws.Range("S32").Value = $s32

# Row 72 "Best Trial" + its fold RMSE list
$ws.Range("K72").Value = "56"
$ws.Range("S72").Value = $s72

# Row 71 fold RMSE list
$ws.Range("S71").Value = $s71

# Row 87 "Best Trial" (reuses "56") + its fold RMSE list
$ws.Range("K87").Value = "56"
$ws.Range("S87").Value = $s87

# --- Row 31 : remaining cells ---
$ws.Range("F31").Value = 0.0391326998151519
$ws.Range("G31").Value = "Adam"
$ws.Range("H31").Value = "relu"
$ws.Range("I31").Value = "1024"
$ws.Range("J31").Value = 0.304653177209175
$ws.Range("K31").Value = "75"
$ws.Range("L31").Value = 0.024652411788702
$ws.Range("P31").Value = 0.00607873736701302

# --- Row 32 : remaining cells ---
$ws.Range("F32").Value = 0.00389624089650994
$ws.Range("G32").Value = "SGD"
$ws.Range("H32").Value = "relu"
$ws.Range("I32").Value = "512"
$ws.Range("J32").Value = 0.160233512133126
$ws.Range("K32").Value = "42"
$ws.Range("L32").Value = 0.0423423804342746
$ws.Range("P32").Value = 0.00847099347178483

# --- Row 47 : remaining cells ---
$ws.Range("F47").Value = 0.0982519160401855
$ws.Range("G47").Value = "Adam"
$ws.Range("H47").Value = "relu"
$ws.Range("I47").Value = "256"
$ws.Range("J47").Value = 0.334652561931301
$ws.Range("L47").Value = 0.0641570977866649
$ws.Range("P47").Value = 0.00713535603902522

# --- Row 71 : remaining cells ---
$ws.Range("E71").Value = "256"
$ws.Range("F71").Value = 0.0010070430958316
$ws.Range("G71").Value = "Adam"
$ws.Range("H71").Value = "tanh"
$ws.Range("I71").Value = "1024"
$ws.Range("J71").Value = 0.28840183199671
$ws.Range("K71").Value = "81"
$ws.Range("L71").Value = 0.0170373452827334
$ws.Range("P71").Value = 0.00383432754987065

# --- Row 72 : remaining cells ---
$ws.Range("E72").Value = "512"
$ws.Range("F72").Value = 0.0273040792261144
$ws.Range("G72").Value = "SGD"
$ws.Range("H72").Value = "tanh"
$ws.Range("I72").Value = "1024"
$ws.Range("J72").Value = 0.318261020392119
$ws.Range("L72").Value = 0.0398057863116264
$ws.Range("P72").Value = 0.00701124318502413

# --- Row 87 : remaining cells ---
$ws.Range("E87").Value = "512"
$ws.Range("F87").Value = 0.00202511691325992
$ws.Range("G87").Value = "Adam"
$ws.Range("H87").Value = "tanh"
$ws.Range("I87").Value = "512"
$ws.Range("J87").Value = 0.140422682361351
$ws.Range("L87").Value = 0.0381521053612232
$ws.Range("P87").Value = 0.00508181241436022

# Recalculate so the dependent Q/R CI formulas pick up real numeric results
$excel.CalculateFullRebuild()
$excel.Calculate()

# --- Update sheet view: scrolled position and selection (matches sheet2.xml sheetView diff) ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$ws.Range("L87").Select()
